# Update cryptocurrency "Price" column (D) values to the latest scrape.
# Values are stored as text (inline strings) in the workbook, so we assign
# them as strings to preserve exact formatting (trailing zeros, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "267.83"
    3  = "21.48"
    4  = "6.245"
    5  = "0.06188"
    6  = "3.565"
    7  = "6.552"
    8  = "1.371"
    9  = "0.8221"
    10 = "0.01347"
    11 = "0.1554"
    12 = "0.08160"
    13 = "0.03299"
    14 = "0.03175"
    15 = "0.09274"
    16 = "3.755"
    17 = "0.001618"
    18 = "0.04681"
    19 = "0.006385"
    20 = "0.005769"
    21 = "0.001068"
    23 = "3.721"
    24 = "2.262"
    40 = "0.04659"
    41 = "0.006972"
    42 = "0.003995"
    44 = "0.01188"
    45 = "0.00006038"
    46 = "0.0009887"
    48 = "0.7810"
    49 = "0.002439"
    50 = "0.00001897"
    51 = "0.01238"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    # Leading apostrophe forces Excel to store the value as literal text,
    # preserving exact formatting (e.g. trailing zeros like "0.7810")
    # instead of silently coercing it to a Double.
    $cell.Value = "'" + $updates[$row]
    # Re-apply the default style so the quote-prefix marker Excel adds
    # doesn't leave a stray cell style behind.
    $cell.Style = "Normal"
}
